# Replace the old manager/password test credentials with the new ones on
# the "Data" sheet (rows 1 and 5 hold the same header/credential pair),
# and move the active selection to E8 (was C11).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$ws.Activate()

$ws.Range("A1").Value = "mngr493079"
$ws.Range("B1").Value = "zeqehEn"
$ws.Range("A5").Value = "mngr493079"
$ws.Range("B5").Value = "zeqehEn"

$ws.Range("E8").Select()
